$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set column B width to 52 characters (calibrated input to hit exact target after Excel's pixel rounding)
$ws.Columns.Item(2).ColumnWidth = 51.15

# Write all row data (A:E, G:H); F (URL) handled separately via Hyperlinks so styling + relationships are correct
$ws.Cells.Item(2,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(2,2).Value = "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,7).Value = 368
$ws.Cells.Item(2,8).Value = "🔥AI,Ai ◆開発"

$ws.Cells.Item(3,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(3,2).Value = "Google AI studio が生成したウェブアプリの補修・ユーザー認証実装"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,7).Value = 338
$ws.Cells.Item(3,8).Value = "🔥AI,Ai ◇アプリ"

$ws.Cells.Item(4,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(4,2).Value = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,7).Value = 310
$ws.Cells.Item(4,8).Value = "🔥AI,Ai"

$ws.Cells.Item(5,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(5,2).Value = "【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,7).Value = 303
$ws.Cells.Item(5,8).Value = "🔥AI,Ai"

$ws.Cells.Item(6,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(6,2).Value = "【急募】メルカリ出品商品の在庫管理自動化ツール開発"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,7).Value = 215
$ws.Cells.Item(6,8).Value = "◆ツール,開発 ◇管理"

$ws.Cells.Item(7,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(7,2).Value = "【Zapier保守・運用サポート】既存フローの管理・調整をお任せできる方募集(時給1,150円程度)"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,7).Value = 213
$ws.Cells.Item(7,8).Value = "🔥API ◇管理"

$ws.Cells.Item(8,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(8,2).Value = "【急募】webアプリ開発のエンジニアを探しています!"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,7).Value = 100
$ws.Cells.Item(8,8).Value = "◆開発 ◇アプリ"

$ws.Cells.Item(9,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(9,2).Value = "【北海道・沖縄】2026年度新人Java研修講師募集!3カ月の短期"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,7).Value = 85
$ws.Cells.Item(9,8).Value = "★Java"

$ws.Cells.Item(10,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(10,2).Value = "【急募】パルワールドのMOD開発に関する依頼"
$ws.Cells.Item(10,3).Value = "システム開発"
$ws.Cells.Item(10,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(10,5).Value = "期限情報なし"
$ws.Cells.Item(10,7).Value = 68
$ws.Cells.Item(10,8).Value = "◆開発"

$ws.Cells.Item(11,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(11,2).Value = "[日本人限定]webシステム制作のプロジェクト管理業務"
$ws.Cells.Item(11,3).Value = "システム開発"
$ws.Cells.Item(11,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(11,5).Value = "期限情報なし"
$ws.Cells.Item(11,7).Value = 60
$ws.Cells.Item(11,8).Value = "◇管理"

$ws.Cells.Item(12,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(12,2).Value = "【急募】WordPressでの会議室予約システム構築依頼"
$ws.Cells.Item(12,3).Value = "システム開発"
$ws.Cells.Item(12,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(12,5).Value = "期限情報なし"
$ws.Cells.Item(12,7).Value = 48
$ws.Cells.Item(12,8).Value = "○WordPress"

$ws.Cells.Item(13,1).Value = "2026-01-27 12:43:44"
$ws.Cells.Item(13,2).Value = "AntigravityからAndroid、iOSを含めたアプリのリリースを教えてほしい"
$ws.Cells.Item(13,3).Value = "システム開発"
$ws.Cells.Item(13,4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(13,5).Value = "期限情報なし"
$ws.Cells.Item(13,7).Value = 30
$ws.Cells.Item(13,8).Value = "◇アプリ"

# Rebuild hyperlinks for column F in row order so relationship IDs line up sequentially.
# TextToDisplay (5th arg) also sets the cell's visible text to the URL, matching the source data.
$ws.Range("F2:F11").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5434128", "", "", "https://www.lancers.jp/work/detail/5434128")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5479941", "", "", "https://www.lancers.jp/work/detail/5479941")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5427956", "", "", "https://www.lancers.jp/work/detail/5427956")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5439158", "", "", "https://www.lancers.jp/work/detail/5439158")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5479836", "", "", "https://www.lancers.jp/work/detail/5479836")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5475245", "", "", "https://www.lancers.jp/work/detail/5475245")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5479608", "", "", "https://www.lancers.jp/work/detail/5479608")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5479693", "", "", "https://www.lancers.jp/work/detail/5479693")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5479655", "", "", "https://www.lancers.jp/work/detail/5479655")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5479860", "", "", "https://www.lancers.jp/work/detail/5479860")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5479809", "", "", "https://www.lancers.jp/work/detail/5479809")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5479715", "", "", "https://www.lancers.jp/work/detail/5479715")
